$d = $word.ActiveDocument

# 1. Update the Sift operations engineer line: replace " (" ... ")" with ": " ... ""
$d.Content.Find.Execute(
    "Operations Engineer, Sift, 2009-2011 (Drupal, CentOS, Xen, VMWare/ESXi)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Operations Engineer, Sift, 2009-2011: Drupal, CentOS, Xen, VMWare/ESXi",
    2
)

# 2. Remove the stray "aaaa" paragraph entirely (including its paragraph mark)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "aaaa") {
        $p.Range.Delete()
        break
    }
}
